# Close #44 - Notifications Controller
# Append three new loan rows (Pepe borrowing twice, Goku borrowing once) to
# the "Loans" sheet, mirroring the existing string-typed "0"/"1" id values
# already used by rows 5-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Loans")

$rows = @(
    @("0", "0", "Pepe", "Wed May 23 02:07:40 ART 2018", "Thu May 24 02:07:40 ART 2018"),
    @("0", "0", "Pepe", "Wed May 23 02:07:42 ART 2018", "Thu May 24 02:07:42 ART 2018"),
    @("0", "1", "Goku", "Wed May 23 02:07:42 ART 2018", "Thu May 24 02:07:42 ART 2018")
)

$startRow = 12
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Columns A & B hold numeric-looking id strings ("0"/"1") that must be
    # stored as TEXT (shared string), matching rows 5-11. Writing the value
    # directly stores a Number, and using a leading apostrophe stores text
    # but tags the cell with a quotePrefix style that the original file
    # doesn't have. Going through a text formula and flattening it with
    # Copy / PasteSpecial (values-only) yields a plain text cell, same as
    # the existing rows.
    $ws.Cells.Item($r, 1).Formula = '="' + $data[0] + '"'
    $ws.Cells.Item($r, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4163)

    $ws.Cells.Item($r, 2).Formula = '="' + $data[1] + '"'
    $ws.Cells.Item($r, 2).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4163)

    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}
